$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H43").Value = 1176.7142
$ws.Range("I43").Value = 799
$ws.Range("J43").Value = 1239.6666
$ws.Range("K43").Value = 799
$ws.Range("L43").Value = 1239.6666
$ws.Range("M43").Value = -730
$ws.Range("N43").Value = -1377.6666
$ws.Range("H51").Value = 2960
$ws.Range("I51").Value = 2900
$ws.Range("K51").Value = 2900
$ws.Range("M51").Value = -2416
$ws.Range("H55").Value = 384.18182
$ws.Range("I55").Value = 46.666668
$ws.Range("J55").Value = 510.75
$ws.Range("K55").Value = 46.666668
$ws.Range("L55").Value = 510.75
$ws.Range("M55").Value = 167.333332
$ws.Range("N55").Value = -938.75
$ws.Range("H132").Value = 2642.375
$ws.Range("I132").Value = 3403.3333
$ws.Range("J132").Value = 359.5
$ws.Range("K132").Value = 10209.9999
$ws.Range("L132").Value = 1078.5
$ws.Range("M132").Value = -7679.999899999999
$ws.Range("N132").Value = -6138.5
$ws.Range("H135").Value = 57693200
$ws.Range("J135").Value = 125002310
$ws.Range("L135").Value = 1125020790
$ws.Range("N135").Value = -1125025860
$ws.Range("H138").Value = 4568.1553
$ws.Range("I138").Value = 2167.6924
$ws.Range("J138").Value = 5261.622
$ws.Range("K138").Value = 6503.0772
$ws.Range("L138").Value = 15784.866
$ws.Range("M138").Value = -1363.0772
$ws.Range("N138").Value = -26064.866
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 70028
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H45").Value = 2082.1304
$ws.Range("I45").Value = 1994.7142
$ws.Range("K45").Value = 1994.7142
$ws.Range("M45").Value = -1617.7142
$ws.Range("H61").Value = 6201.6562
$ws.Range("I61").Value = 4701.1304
$ws.Range("J61").Value = 10036.333
$ws.Range("K61").Value = 4701.1304
$ws.Range("L61").Value = 10036.333
$ws.Range("M61").Value = -4489.1304
$ws.Range("N61").Value = -10460.333
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 1070.0869
$ws.Range("I97").Value = 981.6316
$ws.Range("J97").Value = 1490.25
$ws.Range("K97").Value = 981.6316
$ws.Range("L97").Value = 1490.25
$ws.Range("M97").Value = -485.6316
$ws.Range("N97").Value = -2482.25
$ws.Range("H132").Value = 2709.7441
$ws.Range("I132").Value = 2699.1924
$ws.Range("J132").Value = 2725.8823
$ws.Range("K132").Value = 8097.5772
$ws.Range("L132").Value = 8177.646900000001
$ws.Range("M132").Value = -5567.5772
$ws.Range("N132").Value = -13237.6469
$ws.Range("H136").Value = 6201.6562
$ws.Range("I136").Value = 4701.1304
$ws.Range("J136").Value = 10036.333
$ws.Range("K136").Value = 14103.3912
$ws.Range("L136").Value = 30108.999
$ws.Range("M136").Value = -11553.3912
$ws.Range("N136").Value = -35208.999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1831.8182
$ws.Range("I99").Value = 1215
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 1215
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = 283
$ws.Range("N99").Value = -10996
$ws.Range("H134").Value = 2644
$ws.Range("I134").Value = 2899.4666
$ws.Range("J134").Value = 1366.6666
$ws.Range("K134").Value = 8698.399800000001
$ws.Range("L134").Value = 4099.9998
$ws.Range("M134").Value = -6163.399800000001
$ws.Range("N134").Value = -9169.9998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 444738.22
$ws.Range("I31").Value = 4990.156
$ws.Range("J31").Value = 952139.8
$ws.Range("K31").Value = 4990.156
$ws.Range("L31").Value = 952139.8
$ws.Range("M31").Value = -4695.156
$ws.Range("N31").Value = -952729.8
$ws.Range("H34").Value = 444738.22
$ws.Range("I34").Value = 4990.156
$ws.Range("J34").Value = 952139.8
$ws.Range("K34").Value = 4990.156
$ws.Range("L34").Value = 952139.8
$ws.Range("M34").Value = -4788.156
$ws.Range("N34").Value = -952543.8
$ws.Range("H45").Value = 15000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 15000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 15000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -16186
$ws.Range("H58").Value = 2221366.8
$ws.Range("I58").Value = 3953759
$ws.Range("J58").Value = 7754.4443
$ws.Range("K58").Value = 3953759
$ws.Range("L58").Value = 7754.4443
$ws.Range("M58").Value = -3953556
$ws.Range("N58").Value = -8160.4443
$ws.Range("H134").Value = 2330.3408
$ws.Range("I134").Value = 1902.6333
$ws.Range("J134").Value = 3246.8572
$ws.Range("K134").Value = 5707.8999
$ws.Range("L134").Value = 9740.571599999999
$ws.Range("M134").Value = -3172.8999
$ws.Range("N134").Value = -14810.5716
$ws.Range("H136").Value = 2221366.8
$ws.Range("I136").Value = 3953759
$ws.Range("J136").Value = 7754.4443
$ws.Range("K136").Value = 11861277
$ws.Range("L136").Value = 23263.3329
$ws.Range("M136").Value = -11858727
$ws.Range("N136").Value = -28363.3329
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 35714492
$ws.Range("I12").Value = 125000240
$ws.Range("J12").Value = 192.6
$ws.Range("K12").Value = 375000720
$ws.Range("L12").Value = 577.8
$ws.Range("M12").Value = -375000547
$ws.Range("N12").Value = -923.8
$ws.Range("H42").Value = 3885
$ws.Range("J42").Value = 3885
$ws.Range("L42").Value = 11655
$ws.Range("N42").Value = -12723
$ws.Range("H63").Value = 4866.3335
$ws.Range("J63").Value = 4866.3335
$ws.Range("L63").Value = 14599.0005
$ws.Range("N63").Value = -16097.0005
$ws.Range("H66").Value = 4866.3335
$ws.Range("J66").Value = 4866.3335
$ws.Range("L66").Value = 43797.0015
$ws.Range("N66").Value = -51285.0015
$ws.Range("H75").Value = 3635.7
$ws.Range("J75").Value = 4740.6665
$ws.Range("L75").Value = 14221.9995
$ws.Range("N75").Value = -16217.9995
$ws.Range("H78").Value = 3635.7
$ws.Range("J78").Value = 4740.6665
$ws.Range("L78").Value = 42665.9985
$ws.Range("N78").Value = -52649.9985
$ws.Range("H131").Value = 1297.6522
$ws.Range("I131").Value = 1027.5
$ws.Range("J131").Value = 1441.7333
$ws.Range("K131").Value = 3082.5
$ws.Range("L131").Value = 4325.199900000001
$ws.Range("M131").Value = 1957.5
$ws.Range("N131").Value = -14405.1999
$ws.Range("H137").Value = 24242.13
$ws.Range("I137").Value = 1150.4706
$ws.Range("J137").Value = 89668.5
$ws.Range("K137").Value = 3451.4118
$ws.Range("L137").Value = 269005.5
$ws.Range("M137").Value = 1648.5882
$ws.Range("N137").Value = -279205.5
$ws.Range("H141").Value = 3679.652
$ws.Range("I141").Value = 2201.4285
$ws.Range("K141").Value = 6604.2855
$ws.Range("M141").Value = -1424.2855
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 10031
$ws.Range("J44").Value = 10031
$ws.Range("L44").Value = 10031
$ws.Range("N44").Value = -11223
$ws.Range("H132").Value = 2478.6956
$ws.Range("I132").Value = 2047.1428
$ws.Range("J132").Value = 3150
$ws.Range("K132").Value = 6141.428400000001
$ws.Range("L132").Value = 9450
$ws.Range("M132").Value = -3611.428400000001
$ws.Range("N132").Value = -14510
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4676.25
$ws.Range("I7").Value = 4676.25
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4676.25
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4564.25
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 465.5
$ws.Range("I22").Value = 497.25
$ws.Range("J22").Value = 402
$ws.Range("K22").Value = 497.25
$ws.Range("L22").Value = 402
$ws.Range("M22").Value = -202.25
$ws.Range("N22").Value = -992
$ws.Range("H27").Value = 465.5
$ws.Range("I27").Value = 497.25
$ws.Range("J27").Value = 402
$ws.Range("K27").Value = 497.25
$ws.Range("L27").Value = 402
$ws.Range("M27").Value = -390.25
$ws.Range("N27").Value = -616
$ws.Range("H46").Value = 1160
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("M46").Value = -612
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H93").Value = 867.05884
$ws.Range("I93").Value = 871.25
$ws.Range("K93").Value = 871.25
$ws.Range("M93").Value = 376.75
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H126").Value = 4676.25
$ws.Range("I126").Value = 4676.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14028.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11558.75
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 2474.1167
$ws.Range("I136").Value = 1715.1708
$ws.Range("J136").Value = 4111.8423
$ws.Range("K136").Value = 5145.512400000001
$ws.Range("L136").Value = 12335.5269
$ws.Range("M136").Value = -2595.512400000001
$ws.Range("N136").Value = -17435.5269
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H74").Value = 14252.3
$ws.Range("J74").Value = 14724.777
$ws.Range("L74").Value = 14724.777
$ws.Range("N74").Value = -16596.777
$ws.Range("H77").Value = 14252.3
$ws.Range("J77").Value = 14724.777
$ws.Range("L77").Value = 44174.331
$ws.Range("N77").Value = -53534.331
$ws.Range("H126").Value = 1780.7
$ws.Range("I126").Value = 1742
$ws.Range("K126").Value = 5226
$ws.Range("M126").Value = -2756
